# Re-sort the category rows (A2:C10) alphabetically by the category code in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current rows (2..10) into a list of row-value arrays, keyed by the
# category code in column A so we can reorder them deterministically.
$rowsByCode = @{}
for ($r = 2; $r -le 10; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    $rowsByCode[$code] = @(
        $code,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2
    )
}

# Target order: alphabetical by the K_ code (matches the committed edit).
$orderedCodes = $rowsByCode.Keys | Sort-Object

# Write the sorted rows back into A2:C10
$r = 2
foreach ($code in $orderedCodes) {
    $row = $rowsByCode[$code]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
